# Generate Report for Handoff
#
# The previous handoff round (f42e4ba0-...) has been superseded: its row is
# removed from every sheet, and the still-open handoff (11d03539-...) moves
# from "Handed back: in sync with en-US" to "Ready for handoff" with fresh
# handoff timestamps.

$wb = $excel.ActiveWorkbook

# Helper: repeatedly rescan the Hyperlinks collection and delete the first
# item whose display text matches $pattern, one at a time. Deleting while
# iterating (or indexing with .Item(n)) corrupts this host's Hyperlinks
# collection, so each deletion gets its own fresh foreach scan.
function Remove-MatchingHyperlinks($ws, $pattern) {
    $found = $true
    $guard = 0
    while ($found -and $guard -lt 100) {
        $found = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.TextToDisplay -like $pattern) {
                $hl.Delete()
                $found = $true
                break
            }
        }
        $guard++
    }
}

# ---------------------------------------------------------------------
# Overview sheet: drop the f42e4ba0 row, refresh status + handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows(3).Delete()
Remove-MatchingHyperlinks $wsOverview "f42e4ba0*"

$wsOverview.Cells.Item(2, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 4).Value = "2016-03-30 10:11:57"

# ---------------------------------------------------------------------
# zh-cn sheet: drop the f42e4ba0 row, refresh status + handoff datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows(3).Delete()
Remove-MatchingHyperlinks $wsZhCn "f42e4ba0*"

$wsZhCn.Cells.Item(2, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(2, 5).Value = "2016-03-30 10:11:46"

# ---------------------------------------------------------------------
# de-de sheet: drop the f42e4ba0 row, refresh status + handoff datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows(3).Delete()
Remove-MatchingHyperlinks $wsDeDe "f42e4ba0*"

$wsDeDe.Cells.Item(2, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(2, 5).Value = "2016-03-30 10:11:57"
